# Fruta / hortaliza, semanal
# Insert a new weekly record at row 424 (Terminal La Palmera de La Serena - Ajo),
# shifting the existing historical rows 424:452 down to 425:453.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows down by inserting a fresh row at 424.
$ws.Rows("424:424").Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(424, 1).Value = 8
$ws.Cells.Item(424, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(424, 3).Value = "Coquimbo"
$ws.Cells.Item(424, 4).Value = 45021
$ws.Cells.Item(424, 5).Value = 4
$ws.Cells.Item(424, 6).Value = 100112003
$ws.Cells.Item(424, 7).Value = "Ajo"
$ws.Cells.Item(424, 8).Value = "Chino"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 400
$ws.Cells.Item(424, 11).Value = 17000
$ws.Cells.Item(424, 12).Value = 18000
$ws.Cells.Item(424, 13).Value = 17500
$ws.Cells.Item(424, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(424, 15).Value = "China"
$ws.Cells.Item(424, 16).Value = 1750
$ws.Cells.Item(424, 17).Value = 10
$ws.Cells.Item(424, 18).Value = "Hortaliza"
